$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 7, mirroring the layout/style of the existing data rows.
$ws.Range("A7").Value = 45171
$ws.Range("A7").NumberFormat = "yyyy-mm-dd"

$ws.Range("B7").Value = "23:15"
$ws.Range("C7").Value = "23:15"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "123456789"
$ws.Range("D7").Style = $ws.Range("D2").Style

$ws.Range("E7").Value = "Test"
$ws.Range("F7").Value = "User"
$ws.Range("G7").Value = "test@test.com"
$ws.Range("H7").Value = "Galipatia"
$ws.Range("I7").Value = "Sophomore"
